$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.200778
$ws.Range("H2").Value = 3.602334
$ws.Range("I2").Value = 0.07334464402956961
$ws.Range("J2").Value = 0.08784648530804995
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 54.43165466666667
$ws.Range("N2").Value = 163.294964
$ws.Range("O2").Value = 0.2228930782800698
$ws.Range("P2").Value = 0.2327227899462091
$ws.Range("Q2").Value = 65.36033342733066
$ws.Range("R2").Value = 588.243000845976
$ws.Range("S2").Value = 0.01634801348310671
$ws.Range("T2").Value = 0.02044387914785805
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.200778
$ws.Range("H3").Value = 3.602334
$ws.Range("I3").Value = 0.07334464402956961
$ws.Range("J3").Value = 0.08784648530804995
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 126.7095336666667
$ws.Range("N3").Value = 380.128601
$ws.Range("O3").Value = 0.5188649542136915
$ws.Range("P3").Value = 0.541747194133123
$ws.Range("Q3").Value = 152.1500204171926
$ws.Range("R3").Value = 1369.350183754734
$ws.Range("S3").Value = 0.03805596536622214
$ws.Range("T3").Value = 0.04759058693009267
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.200778
$ws.Range("H4").Value = 3.602334
$ws.Range("I4").Value = 0.07334464402956961
$ws.Range("J4").Value = 0.08784648530804995
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 17.491284
$ws.Range("N4").Value = 52.473852
$ws.Range("O4").Value = 0.07162534664261168
$ws.Range("P4").Value = 0.07478406521259567
$ws.Range("Q4").Value = 21.003149018952
$ws.Range("R4").Value = 189.028341170568
$ws.Range("S4").Value = 0.005253335552996882
$ws.Range("T4").Value = 0.006569517285974535
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.200778
$ws.Range("H5").Value = 3.602334
$ws.Range("I5").Value = 0.07334464402956961
$ws.Range("J5").Value = 0.08784648530804995
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 14.62863866666667
$ws.Range("N5").Value = 43.88591599999999
$ws.Range("O5").Value = 0.05990305316690945
$ws.Range("P5").Value = 0.06254481192001105
$ws.Range("Q5").Value = 17.56574748088266
$ws.Range("R5").Value = 158.091727327944
$ws.Range("S5").Value = 0.004393568110811356
$ws.Range("T5").Value = 0.005494341901425998
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.200778
$ws.Range("H6").Value = 3.602334
$ws.Range("I6").Value = 0.07334464402956961
$ws.Range("J6").Value = 0.08784648530804995
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 30.9441155
$ws.Range("N6").Value = 61.888231
$ws.Range("O6").Value = 0.1267135676967176
$ws.Range("P6").Value = 0.08820113878806125
$ws.Range("Q6").Value = 37.157013121859
$ws.Range("R6").Value = 222.942078731154
$ws.Range("S6").Value = 0.009293761516432521
$ws.Range("T6").Value = 0.007748160042698697
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.197813
$ws.Range("H7").Value = 9.593439
$ws.Range("I7").Value = 0.1953254108237577
$ws.Range("J7").Value = 0.2339455192570077
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 54.43165466666667
$ws.Range("N7").Value = 163.294964
$ws.Range("O7").Value = 0.2228930782800698
$ws.Range("P7").Value = 0.2327227899462091
$ws.Range("Q7").Value = 174.0622529045773
$ws.Range("R7").Value = 1566.560276141196
$ws.Range("S7").Value = 0.0435366820848266
$ws.Range("T7").Value = 0.0544444539369054
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.197813
$ws.Range("H8").Value = 9.593439
$ws.Range("I8").Value = 0.1953254108237577
$ws.Range("J8").Value = 0.2339455192570077
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 126.7095336666667
$ws.Range("N8").Value = 380.128601
$ws.Range("O8").Value = 0.5188649542136915
$ws.Range("P8").Value = 0.541747194133123
$ws.Range("Q8").Value = 405.1933939832044
$ws.Range("R8").Value = 3646.740545848839
$ws.Range("S8").Value = 0.1013475103438395
$ws.Range("T8").Value = 0.1267393286375004
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.197813
$ws.Range("H9").Value = 9.593439
$ws.Range("I9").Value = 0.1953254108237577
$ws.Range("J9").Value = 0.2339455192570077
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.491284
$ws.Range("N9").Value = 52.473852
$ws.Range("O9").Value = 0.07162534664261168
$ws.Range("P9").Value = 0.07478406521259567
$ws.Range("Q9").Value = 55.933855361892
$ws.Range("R9").Value = 503.404698257028
$ws.Range("S9").Value = 0.01399025025836218
$ws.Range("T9").Value = 0.01749539696831062
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.197813
$ws.Range("H10").Value = 9.593439
$ws.Range("I10").Value = 0.1953254108237577
$ws.Range("J10").Value = 0.2339455192570077
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.62863866666667
$ws.Range("N10").Value = 43.88591599999999
$ws.Range("O10").Value = 0.05990305316690945
$ws.Range("P10").Value = 0.06254481192001105
$ws.Range("Q10").Value = 46.77965090056933
$ws.Range("R10").Value = 421.016858105124
$ws.Range("S10").Value = 0.01170058846942399
$ws.Range("T10").Value = 0.01463207850145887
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 3.197813
$ws.Range("H11").Value = 9.593439
$ws.Range("I11").Value = 0.1953254108237577
$ws.Range("J11").Value = 0.2339455192570077
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 30.9441155
$ws.Range("N11").Value = 61.888231
$ws.Range("O11").Value = 0.1267135676967176
$ws.Range("P11").Value = 0.08820113878806125
$ws.Range("Q11").Value = 98.95349481940151
$ws.Range("R11").Value = 593.720968916409
$ws.Range("S11").Value = 0.02475037966730539
$ws.Range("T11").Value = 0.02063426121283239
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.622028
$ws.Range("H12").Value = 4.866084
$ws.Range("I12").Value = 0.09907498827093329
$ws.Range("J12").Value = 0.118664281716725
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.43165466666667
$ws.Range("N12").Value = 163.294964
$ws.Range("O12").Value = 0.2228930782800698
$ws.Range("P12").Value = 0.2327227899462091
$ws.Range("Q12").Value = 88.289667955664
$ws.Range("R12").Value = 794.607011600976
$ws.Range("S12").Value = 0.02208312911627013
$ws.Range("T12").Value = 0.02761588270807918
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.622028
$ws.Range("H13").Value = 4.866084
$ws.Range("I13").Value = 0.09907498827093329
$ws.Range("J13").Value = 0.118664281716725
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 126.7095336666667
$ws.Range("N13").Value = 380.128601
$ws.Range("O13").Value = 0.5188649542136915
$ws.Range("P13").Value = 0.541747194133123
$ws.Range("Q13").Value = 205.526411474276
$ws.Range("R13").Value = 1849.737703268484
$ws.Range("S13").Value = 0.05140653925291983
$ws.Range("T13").Value = 0.06428604166385823
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.622028
$ws.Range("H14").Value = 4.866084
$ws.Range("I14").Value = 0.09907498827093329
$ws.Range("J14").Value = 0.118664281716725
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 17.491284
$ws.Range("N14").Value = 52.473852
$ws.Range("O14").Value = 0.07162534664261168
$ws.Range("P14").Value = 0.07478406521259567
$ws.Range("Q14").Value = 28.371352403952
$ws.Range("R14").Value = 255.342171635568
$ws.Range("S14").Value = 0.007096280378518284
$ws.Range("T14").Value = 0.008874197382309389
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.622028
$ws.Range("H15").Value = 4.866084
$ws.Range("I15").Value = 0.09907498827093329
$ws.Range("J15").Value = 0.118664281716725
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 14.62863866666667
$ws.Range("N15").Value = 43.88591599999999
$ws.Range("O15").Value = 0.05990305316690945
$ws.Range("P15").Value = 0.06254481192001105
$ws.Range("Q15").Value = 23.728061519216
$ws.Range("R15").Value = 213.552553672944
$ws.Range("S15").Value = 0.005934894289904647
$ws.Range("T15").Value = 0.007421835181595772
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.622028
$ws.Range("H16").Value = 4.866084
$ws.Range("I16").Value = 0.09907498827093329
$ws.Range("J16").Value = 0.118664281716725
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 30.9441155
$ws.Range("N16").Value = 61.888231
$ws.Range("O16").Value = 0.1267135676967176
$ws.Range("P16").Value = 0.08820113878806125
$ws.Range("Q16").Value = 50.19222177623401
$ws.Range("R16").Value = 301.153330657404
$ws.Range("S16").Value = 0.01255414523332041
$ws.Range("T16").Value = 0.01046632478088246
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.243090333333333
$ws.Range("H17").Value = 6.729271
$ws.Range("I17").Value = 0.1370100568335712
$ws.Range("J17").Value = 0.1640999435464303
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 54.43165466666667
$ws.Range("N17").Value = 163.294964
$ws.Range("O17").Value = 0.2228930782800698
$ws.Range("P17").Value = 0.2327227899462091
$ws.Range("Q17").Value = 122.0951184101382
$ws.Range("R17").Value = 1098.856065691244
$ws.Range("S17").Value = 0.030538593322962
$ws.Range("T17").Value = 0.03818979669214067
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 2.243090333333333
$ws.Range("H18").Value = 6.729271
$ws.Range("I18").Value = 0.1370100568335712
$ws.Range("J18").Value = 0.1640999435464303
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 126.7095336666667
$ws.Range("N18").Value = 380.128601
$ws.Range("O18").Value = 0.5188649542136915
$ws.Range("P18").Value = 0.541747194133123
$ws.Range("Q18").Value = 284.2209301088746
$ws.Range("R18").Value = 2557.988370979871
$ws.Range("S18").Value = 0.07108971686576621
$ws.Range("T18").Value = 0.0889006839736825
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 2.243090333333333
$ws.Range("H19").Value = 6.729271
$ws.Range("I19").Value = 0.1370100568335712
$ws.Range("J19").Value = 0.1640999435464303
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 17.491284
$ws.Range("N19").Value = 52.473852
$ws.Range("O19").Value = 0.07162534664261168
$ws.Range("P19").Value = 0.07478406521259567
$ws.Range("Q19").Value = 39.234530057988
$ws.Range("R19").Value = 353.110770521892
$ws.Range("S19").Value = 0.009813392814228467
$ws.Range("T19").Value = 0.01227206087955951
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 2.243090333333333
$ws.Range("H20").Value = 6.729271
$ws.Range("I20").Value = 0.1370100568335712
$ws.Range("J20").Value = 0.1640999435464303
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 14.62863866666667
$ws.Range("N20").Value = 43.88591599999999
$ws.Range("O20").Value = 0.05990305316690945
$ws.Range("P20").Value = 0.06254481192001105
$ws.Range("Q20").Value = 32.81335798302622
$ws.Range("R20").Value = 295.3202218472359
$ws.Range("S20").Value = 0.008207320718902704
$ws.Range("T20").Value = 0.01026360010519591
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 2.243090333333333
$ws.Range("H21").Value = 6.729271
$ws.Range("I21").Value = 0.1370100568335712
$ws.Range("J21").Value = 0.1640999435464303
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 30.9441155
$ws.Range("N21").Value = 61.888231
$ws.Range("O21").Value = 0.1267135676967176
$ws.Range("P21").Value = 0.08820113878806125
$ws.Range("Q21").Value = 69.41044635160017
$ws.Range("R21").Value = 416.462678109601
$ws.Range("S21").Value = 0.01736103311171185
$ws.Range("T21").Value = 0.01447380189585172
$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 8.108010999999999
$ws.Range("H22").Value = 16.216022
$ws.Range("I22").Value = 0.4952449000421683
$ws.Range("J22").Value = 0.3954437701717871
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 54.43165466666667
$ws.Range("N22").Value = 163.294964
$ws.Range("O22").Value = 0.2228930782800698
$ws.Range("P22").Value = 0.2327227899462091
$ws.Range("Q22").Value = 441.3324547855347
$ws.Range("R22").Value = 2647.994728713208
$ws.Range("S22").Value = 0.1103866602729043
$ws.Range("T22").Value = 0.09202877746122579
$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 8.108010999999999
$ws.Range("H23").Value = 16.216022
$ws.Range("I23").Value = 0.4952449000421683
$ws.Range("J23").Value = 0.3954437701717871
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 126.7095336666667
$ws.Range("N23").Value = 380.128601
$ws.Range("O23").Value = 0.5188649542136915
$ws.Range("P23").Value = 0.541747194133123
$ws.Range("Q23").Value = 1027.362292774204
$ws.Range("R23").Value = 6164.173756645222
$ws.Range("S23").Value = 0.2569652223849438
$ws.Range("T23").Value = 0.2142305529279892
$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 8.108010999999999
$ws.Range("H24").Value = 16.216022
$ws.Range("I24").Value = 0.4952449000421683
$ws.Range("J24").Value = 0.3954437701717871
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 17.491284
$ws.Range("N24").Value = 52.473852
$ws.Range("O24").Value = 0.07162534664261168
$ws.Range("P24").Value = 0.07478406521259567
$ws.Range("Q24").Value = 141.819523076124
$ws.Range("R24").Value = 850.9171384567439
$ws.Range("S24").Value = 0.03547208763850587
$ws.Range("T24").Value = 0.02957289269644162
$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 8.108010999999999
$ws.Range("H25").Value = 16.216022
$ws.Range("I25").Value = 0.4952449000421683
$ws.Range("J25").Value = 0.3954437701717871
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 14.62863866666667
$ws.Range("N25").Value = 43.88591599999999
$ws.Range("O25").Value = 0.05990305316690945
$ws.Range("P25").Value = 0.06254481192001105
$ws.Range("Q25").Value = 118.6091632243587
$ws.Range("R25").Value = 711.6549793461519
$ws.Range("S25").Value = 0.02966668157786676
$ws.Range("T25").Value = 0.0247329562303345
$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 8.108010999999999
$ws.Range("H26").Value = 16.216022
$ws.Range("I26").Value = 0.4952449000421683
$ws.Range("J26").Value = 0.3954437701717871
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 30.9441155
$ws.Range("N26").Value = 61.888231
$ws.Range("O26").Value = 0.1267135676967176
$ws.Range("P26").Value = 0.08820113878806125
$ws.Range("Q26").Value = 250.8952288592705
$ws.Range("R26").Value = 1003.580915437082
$ws.Range("S26").Value = 0.06275424816794742
$ws.Range("T26").Value = 0.03487859085579599
